$d = $word.ActiveDocument

# --- Step 1: remove the old "_GoBack" bookmark from the "Datum" paragraph
#     (it currently sits between "...-06" and "-2013").
$oldBm = $d.Bookmarks("_GoBack")
$oldBm.Delete()

# --- Step 2: replace the placeholder "..." paragraph with the real
#     minutes text, conceptually split (like the source edit) into three
#     segments:
#       s1  ->  "We hebben eerst ... plaatsvind "
#       s2  ->  "ook op deze manier"
#       s3  ->  " moeten toelichten. ... relaties."
#     with the "_GoBack" bookmark re-inserted between s2 and s3 (this is
#     where the author's cursor/last edit ended up).

$s1 = "We hebben eerst verteld hoe het gaat en hoever we met de voortgang zijn. We waren blij te kunnen vertellen dat het nog steeds volgens planning verloopt.  Daarna hebben we het over de samenwerking gehad met de andere groep, we hebben uitgelegd dat dit i.v.m. o.a. andere structuur erg lastig wordt, dit begreep onze coach wel maar zei dat we dit straks tijdens het volgende gesprek dat vanmiddag plaatsvind "
$s2 = "ook op deze manier"
$s3 = " moeten toelichten.  We hebben verteld wat er nog moet gebeuren (onze to do list laten zien).  Verder hebben we vermeld dat we nog geen server informatie hebben, welke de andere groep achter aan zou gaan. Tot slot hebben we over de database structuur gesproken wat betreft de vakken, toetsen en inschrijven en de bijbehorende relaties."

# Locate the placeholder paragraph ("...") - it is the last paragraph in the body.
$count = $d.Paragraphs.Count
$target = $d.Paragraphs($count).Range

# Replace its text with the first segment.
$found = $target.Find.Execute("…", $true, $false, $false, $false, $false, $true, 1, $false, $s1, 2)

# Append the second segment right after the first.
$full = $d.Content.Text
$insertPos = $full.IndexOf($s1) + $s1.Length
$afterS1 = $d.Range($insertPos, $insertPos)
$afterS1.InsertAfter($s2)

# Append the third (final) segment right after the second, so the bookmark
# position we compute next is no longer the very last character of the
# document (inserting the bookmark there first, before any trailing text
# exists, lands it in the wrong place).
$full2 = $d.Content.Text
$bmPos = $full2.IndexOf($s2) + $s2.Length
$tail = $d.Range($bmPos, $bmPos)
$tail.InsertAfter($s3)

# Now (re-)locate the boundary between s2 and s3 and drop the "_GoBack"
# bookmark there - it's now a genuine interior position.
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
